# "added 4wk low sales check"
# - Clear the "Inventory Coverage" (H) column for rows 2-14 on "Forecast Comparison"
#   (stockout/low-sales rows no longer report a numeric coverage value).
# - Recompute "Seasonality Index" (L) column values on "Forecast Comparison" for rows 2-17.
# - Zero out the short-horizon forecast totals (4/8/16 week) on "Summary" since the
#   4-week-low-sales check now suppresses them.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Clear Inventory Coverage (column H) for rows 2 through 14
for ($row = 2; $row -le 14; $row++) {
    $wsForecast.Cells.Item($row, 8).Value = $null
}

# Updated Seasonality Index (column L) values for rows 2 through 17
$seasonalityIndex = @{
    2  = 1.03
    3  = 0.8100000000000001
    4  = 0.9399999999999999
    5  = 0.95
    6  = 1.12
    7  = 0.86
    8  = 0.83
    9  = 0.91
    10 = 0.96
    11 = 1.12
    12 = 1.19
    13 = 0.96
    15 = 1.08
    16 = 1.03
    17 = 0.98
}

foreach ($row in $seasonalityIndex.Keys) {
    $wsForecast.Cells.Item($row, 12).Value = $seasonalityIndex[$row]
}

# Zero out short-horizon forecast totals on the Summary sheet.
# These cells hold numbers-as-text (matching the rest of column B, e.g. B3/B12/B14),
# so a leading apostrophe forces the literal text "0" rather than the number 0.
$wsSummary.Range("B9").Value = "'0"
$wsSummary.Range("B10").Value = "'0"
$wsSummary.Range("B11").Value = "'0"
